# Update Name of Algo
# Apply updated imputed values for columns A and B on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B11"  = 6.180000000000001
    "A12"  = -21.37859999999999
    "B23"  = 8.901999999999997
    "A27"  = -21.95550000000001
    "B28"  = 5.656100000000002
    "A32"  = -21.06409999999999
    "B32"  = 6.108699999999998
    "B34"  = 9.713900000000006
    "A36"  = -19.9724
    "A38"  = -20.04339999999998
    "B42"  = 10.31409999999999
    "A46"  = -21.8237
    "B49"  = 5.4788
    "A54"  = -21.684
    "B54"  = 4.512499999999998
    "A55"  = -22.28160000000001
    "A56"  = -22.14190000000001
    "A67"  = -21.38619999999997
    "A69"  = -21.55789999999998
    "A72"  = -21.9857
    "B78"  = 9.663999999999996
    "B80"  = 9.656100000000002
    "A83"  = -21.64199999999999
    "A86"  = -21.933
    "A91"  = -20.9837
    "A93"  = -21.4734
    "B97"  = 5.701199999999997
    "A99"  = -21.78839999999999
    "B99"  = 6.086999999999997
    "B101" = 4.089599999999999
    "A104" = -21.58319999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
